# Auto-generated edit script: update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.552.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.06%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.217.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.44%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "391.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.218.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.95%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.561"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.33%  "

# Row 9
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.611"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "38.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0953"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.54%  "

# Row 13
$ws.Range("E13").Value = "  +1.67%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.741.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.95%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.234.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.77%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.77%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "56.581.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000104"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.85%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.56%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "295.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.76%  "

# Row 27
$ws.Range("E27").Value = "  +2.89%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.68%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.168"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.99%  "

# Row 31
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0479"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.27%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.74%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.56%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.58%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "133.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.89%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.118"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.52%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.282"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.59%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.78%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.39%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.118.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.40%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.16%  "
